{"js": "// Resize the diagram picture in the document:\n//   wp:extent        5731510 x 3415030  ->  5467264 x 3257583\n//   wp:effectExtent  b=1270             ->  b=0\n//   a:ext (shape xfrm) 5731510 x 3415030 -> 5499069 x 3276534\n//   wp14:editId      4DDDF5DD           ->  0F3E8FAC\n//\n// The outer \"display\" extent (wp:extent) and the inner shape extent\n// (a:ext inside pic:spPr/a:xfrm) end up with two DIFFERENT sizes, which\n// can't be produced by setting InlinePicture.width/height alone (that API\n// keeps both in lock-step). So we pull the picture's OOXML, patch the\n// handful of numeric attributes precisely, and write it back in place.\n\nconst pics = context.document.body.inlinePictures;\npics.load(\"items\");\nawait context.sync();\n\nif (pics.items.length === 0) {\n  throw new Error(\"No inline picture found in the document body.\");\n}\n\nconst pic = pics.items[0];\nconst range = pic.getRange();\nconst ooxmlResult = range.getOoxml();\nawait context.sync();\n\nlet xml = ooxmlResult.value;\n\nconst replacements = [\n  ['<wp:extent cx=\"5731510\" cy=\"3415030\"/>', '<wp:extent cx=\"5467264\" cy=\"3257583\"/>'],\n  ['<wp:effectExtent l=\"0\" t=\"0\" r=\"0\" b=\"1270\"/>', '<wp:effectExtent l=\"0\" t=\"0\" r=\"0\" b=\"0\"/>'],\n  ['<a:ext cx=\"5731510\" cy=\"3415030\"/>', '<a:ext cx=\"5499069\" cy=\"3276534\"/>'],\n  ['wp14:editId=\"4DDDF5DD\"', 'wp14:editId=\"0F3E8FAC\"'],\n];\n\nfor (const [from, to] of replacements) {\n  if (xml.indexOf(from) === -1) {\n    throw new Error(\"Expected picture OOXML fragment not found: \" + from);\n  }\n  xml = xml.split(from).join(to);\n}\n\nrange.insertOoxml(xml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Resize the diagram picture in the document:\n#   wp:extent        5731510 x 3415030  ->  5467264 x 3257583\n#   wp:effectExtent  b=1270             ->  b=0\n#   a:ext (shape xfrm) 5731510 x 3415030 -> 5499069 x 3276534\n#   wp14:editId      4DDDF5DD           ->  0F3E8FAC\n#\n# The outer \"display\" extent (wp:extent) and the inner shape extent\n# (a:ext inside pic:spPr/a:xfrm) end up with two DIFFERENT sizes, which\n# InlineShape.Width/Height can't do alone (they keep both in lock-step).\n# So grab the picture's WordOpenXML, patch the numeric attributes\n# precisely, and reinsert it in place of the original range.\n\n$d = $word.ActiveDocument\n\nif ($d.InlineShapes.Count -lt 1) {\n    throw \"No inline picture found in the document body.\"\n}\n\n$shp = $d.InlineShapes.Item(1)\n$rng = $shp.Range\n$xml = $rng.WordOpenXML\n\n$xml = $xml.Replace('<wp:extent cx=\"5731510\" cy=\"3415030\"/>', '<wp:extent cx=\"5467264\" cy=\"3257583\"/>')\n$xml = $xml.Replace('<wp:effectExtent l=\"0\" t=\"0\" r=\"0\" b=\"1270\"/>', '<wp:effectExtent l=\"0\" t=\"0\" r=\"0\" b=\"0\"/>')\n$xml = $xml.Replace('<a:ext cx=\"5731510\" cy=\"3415030\"/>', '<a:ext cx=\"5499069\" cy=\"3276534\"/>')\n$xml = $xml.Replace('wp14:editId=\"4DDDF5DD\"', 'wp14:editId=\"0F3E8FAC\"')\n\n# Delete the original picture range first, then insert the patched XML at\n# the same spot -- InsertXML on a still-populated range appends after it\n# instead of replacing its contents.\n$rng.Delete()\n$rng.InsertXML($xml)\n"}
